$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timing data for rows 2..101 (columns A=Heap, B=Quick, C=Median),
# replacing the benchmark measurements per the HeapSelect O(n)+O(k log k) change.
$data = @(
    @(35127459900,7488199,1754300),
    @(32501,1599,5100),
    @(27300,5850,2650),
    @(13300,1916,5450),
    @(22800,1733,5100),
    @(25899,1900,5249),
    @(20200,1750,10201),
    @(21500,2099,4033),
    @(22000,2900,5649),
    @(18601,1766,4067),
    @(18200,2080,5850),
    @(17999,21099,5550),
    @(16200,1571,3699),
    @(19200,1599,5549),
    @(15800,3766,6600),
    @(50301,1885,3799),
    @(15100,2000,5300),
    @(35001,1300,27999),
    @(17301,1866,4433),
    @(16600,1733,4166),
    @(16500,1499,28900),
    @(99601,1700,5100),
    @(16999,1571,8200),
    @(19401,1816,4133),
    @(20000,1683,5900),
    @(43200,1866,3900),
    @(36199,2383,3600),
    @(16300,1557,5200),
    @(20200,10199,3466),
    @(11100,4100,3433),
    @(16600,1800,3733),
    @(30700,6200,1850),
    @(25701,7500,3866),
    @(20500,1733,8700),
    @(10899,900,2200),
    @(11099,927,1799),
    @(11600,1443,2060),
    @(19000,6450,2119),
    @(66200,6274,3833),
    @(19001,656,3800),
    @(14600,5399,7333),
    @(16800,351,5450),
    @(17100,308,3933),
    @(11800,2308120,3699),
    @(13399,271,3633),
    @(12600,283,12500),
    @(12000,297,10299),
    @(11999,229,5750),
    @(6399,190,12101),
    @(6450,364,1020),
    @(5950,160,1133),
    @(6550,367,2440),
    @(5949,159,1000),
    @(5749,353,963),
    @(18399,534,5150),
    @(12400,294,1883),
    @(13000,400,1766),
    @(12701,305,1337),
    @(14800,490,1374),
    @(13501,329,1166),
    @(12700,302,1049),
    @(12500,315,1090),
    @(11700,283,1040),
    @(11799,300,4450),
    @(13299,312,1030),
    @(61100,329,1287),
    @(13200,447,1716),
    @(28900,294,1110),
    @(28101,309,1059),
    @(10899,360,1059),
    @(13400,265,963),
    @(8100,286,1716),
    @(7800,2058,954),
    @(12000,259,1144),
    @(11300,351,1020),
    @(10300,281,700),
    @(5600,336,600),
    @(5099,192,1059),
    @(7300,542,572),
    @(5500,194,1155),
    @(7300,156,510),
    @(5250,158,520),
    @(4500,156,510),
    @(18099,477,271),
    @(11001,318,443),
    @(10700,312,552),
    @(25200,420,536),
    @(10499,261,519),
    @(6449,161,728),
    @(8300,242,536),
    @(5199,583,275),
    @(10300,200,271),
    @(15100,190,273),
    @(4166,180,268),
    @(4133,1161,374),
    @(4199,271,460),
    @(7349,278,377),
    @(11001,306,536),
    @(4200,429,815),
    @(16800,201,283)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $ws.Cells.Item($i + 2, 1).Value = $row[0]
    $ws.Cells.Item($i + 2, 2).Value = $row[1]
    $ws.Cells.Item($i + 2, 3).Value = $row[2]
}
